# Weekly update: a new daily price record was inserted above the existing
# row 238 ("Pepino ensalada" subset), pushing the prior rows 238-257 down to
# 239-258.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 238; Excel shifts rows 238:257
# down to 239:258 and the sheet's used range grows to A1:R258.
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row 238 with the new record's data.
$ws.Range("A238").Value = 7
$ws.Range("B238").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C238").Value = "Ñuble"
$ws.Range("D238").Value = 44826
$ws.Range("E238").Value = 16
$ws.Range("F238").Value = 100112043
$ws.Range("G238").Value = "Pepino ensalada"
$ws.Range("H238").Value = "Sin especificar"
$ws.Range("I238").Value = "Primera"
$ws.Range("J238").Value = 120
$ws.Range("K238").Value = 19000
$ws.Range("L238").Value = 20000
$ws.Range("M238").Value = 19500
$ws.Range("N238").Value = "$/caja 60 unidades"
$ws.Range("O238").Value = "Región de Arica y Parinacota"
$ws.Range("P238").Value = 325
$ws.Range("Q238").Value = 60
$ws.Range("R238").Value = "Hortaliza"
